$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2962962962962963
$ws.Range("C2").Value = 0.4074074074074074
$ws.Range("J2").Value = 0.07407407407407407
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("J4").Value = 0.1666666666666667
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.1071428571428571
$ws.Range("F6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.3214285714285715
$ws.Range("Q6").Value = 0.2142857142857143
$ws.Range("S6").Value = 0.2857142857142857
$ws.Range("B7").Value = 0.1
$ws.Range("Q7").Value = 0.2
$ws.Range("S7").Value = 0.5
$ws.Range("B8").Value = 0.09302325581395349
$ws.Range("F8").Value = 0.1162790697674419
$ws.Range("J8").Value = 0.06976744186046512
$ws.Range("O8").Value = 0.02325581395348837
$ws.Range("Q8").Value = 0.3720930232558139
$ws.Range("R8").Value = 0.1162790697674419
$ws.Range("S8").Value = 0.2093023255813954
$ws.Range("B9").Value = 0.04
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.04
$ws.Range("Q9").Value = 0.28
$ws.Range("R9").Value = 0.12
$ws.Range("S9").Value = 0.44
$ws.Range("B10").Value = 0.06097560975609756
$ws.Range("D10").Value = 0.03658536585365853
$ws.Range("F10").Value = 0.04878048780487805
$ws.Range("J10").Value = 0.1585365853658537
$ws.Range("O10").Value = 0.01219512195121951
$ws.Range("Q10").Value = 0.25
$ws.Range("R10").Value = 0.1219512195121951
$ws.Range("S10").Value = 0.3109756097560976
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.1333333333333333
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4666666666666667
$ws.Range("G12").Value = 0.8571428571428571
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.04761904761904762
$ws.Range("I15").Value = 0.09523809523809523
$ws.Range("J15").Value = 0.5238095238095238
$ws.Range("K15").Value = 0.04761904761904762
$ws.Range("O15").Value = 0.1904761904761905
$ws.Range("S15").Value = 0.09523809523809523
$ws.Range("F16").Value = 0.05882352941176471
$ws.Range("H16").Value = 0.2352941176470588
$ws.Range("I16").Value = 0.05882352941176471
$ws.Range("J16").Value = 0.4117647058823529
$ws.Range("K16").Value = 0.05882352941176471
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.1176470588235294
$ws.Range("F17").Value = 0.02777777777777778
$ws.Range("H17").Value = 0.2222222222222222
$ws.Range("I17").Value = 0.06944444444444445
$ws.Range("J17").Value = 0.3888888888888889
$ws.Range("K17").Value = 0.1111111111111111
$ws.Range("M17").Value = 0.02777777777777778
$ws.Range("O17").Value = 0.06944444444444445
$ws.Range("S17").Value = 0.08333333333333333
$ws.Range("F18").Value = 0.03703703703703703
$ws.Range("H18").Value = 0.07407407407407407
$ws.Range("I18").Value = 0.2592592592592592
$ws.Range("J18").Value = 0.5555555555555556
$ws.Range("O18").Value = 0.07407407407407407
$ws.Range("F19").Value = 0.0198019801980198
$ws.Range("H19").Value = 0.2079207920792079
$ws.Range("I19").Value = 0.09900990099009901
$ws.Range("J19").Value = 0.5445544554455446
$ws.Range("K19").Value = 0.009900990099009901
$ws.Range("N19").Value = 0.009900990099009901
$ws.Range("O19").Value = 0.04950495049504951
$ws.Range("S19").Value = 0.0594059405940594
